# "Generate Report for handback"
#
# The handback round-trip completed and came back identical to en-US, so:
#  - every "Ready for handoff" status cell becomes
#    "Handed back: in sync with en-US" (Overview sheet + both language sheets)
#  - each language sheet gains a "Latest Target File" (E) / "Latest Handback
#    File" (F) hyperlink pair per row, mirroring the existing handoff file
#    name/link in columns A/C
#  - the "Latest Handback DateTime" (G) placeholder timestamps are stamped
#    with the real handback time

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Hyperlink colour used throughout the workbook's "HyperLink" style
# (font color FF6495ED == RGB(100,149,237), stored as BGR for the COM API).
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet: both language-status columns move to "handed back"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/6cdc25e2dc3afad128f129b27cb8c6d94efba0c8/e2e/0db6ae1f-c28e-4121-a1e7-6223b5cecd9a.md", "", "", "0db6ae1f-c28e-4121-a1e7-6223b5cecd9a.md")
Style-AsHyperlink $wsZh.Range("E2")

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a2d26533643f7b74c80577b252561caac0979e02/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0db6ae1f-c28e-4121-a1e7-6223b5cecd9a.36a13068b54ca1e053de320aef40abcf972cc5b1.zh-cn.xlf", "", "", "0db6ae1f-c28e-4121-a1e7-6223b5cecd9a.36a13068b54ca1e053de320aef40abcf972cc5b1.zh-cn.xlf")
Style-AsHyperlink $wsZh.Range("F2")

$wsZh.Range("G2").Value = "2016-01-28 11:01:00"

$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/6cdc25e2dc3afad128f129b27cb8c6d94efba0c8/e2e/7234652a-7657-405c-9d1c-0e201f8360ac.md", "", "", "7234652a-7657-405c-9d1c-0e201f8360ac.md")
Style-AsHyperlink $wsZh.Range("E3")

$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a2d26533643f7b74c80577b252561caac0979e02/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7234652a-7657-405c-9d1c-0e201f8360ac.0168f3b46f092590a8b80345e5c57900c6671c30.zh-cn.xlf", "", "", "7234652a-7657-405c-9d1c-0e201f8360ac.0168f3b46f092590a8b80345e5c57900c6671c30.zh-cn.xlf")
Style-AsHyperlink $wsZh.Range("F3")

$wsZh.Range("G3").Value = "2016-01-28 11:01:00"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/6cdc25e2dc3afad128f129b27cb8c6d94efba0c8/e2e/0db6ae1f-c28e-4121-a1e7-6223b5cecd9a.md", "", "", "0db6ae1f-c28e-4121-a1e7-6223b5cecd9a.md")
Style-AsHyperlink $wsDe.Range("E2")

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/840b7f66e9a818be686a825c7e53e7a15d544862/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0db6ae1f-c28e-4121-a1e7-6223b5cecd9a.36a13068b54ca1e053de320aef40abcf972cc5b1.de-de.xlf", "", "", "0db6ae1f-c28e-4121-a1e7-6223b5cecd9a.36a13068b54ca1e053de320aef40abcf972cc5b1.de-de.xlf")
Style-AsHyperlink $wsDe.Range("F2")

$wsDe.Range("G2").Value = "2016-01-28 11:01:33"

$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/6cdc25e2dc3afad128f129b27cb8c6d94efba0c8/e2e/7234652a-7657-405c-9d1c-0e201f8360ac.md", "", "", "7234652a-7657-405c-9d1c-0e201f8360ac.md")
Style-AsHyperlink $wsDe.Range("E3")

$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/840b7f66e9a818be686a825c7e53e7a15d544862/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7234652a-7657-405c-9d1c-0e201f8360ac.0168f3b46f092590a8b80345e5c57900c6671c30.de-de.xlf", "", "", "7234652a-7657-405c-9d1c-0e201f8360ac.0168f3b46f092590a8b80345e5c57900c6671c30.de-de.xlf")
Style-AsHyperlink $wsDe.Range("F3")

$wsDe.Range("G3").Value = "2016-01-28 11:01:33"
